$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "67.043.65"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "3.216.12"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "604.75"
$ws.Range("E5").Value = "  +4.40%  "
$ws.Range("D6").Value = "157.56"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +6.32%  "
$ws.Range("D9").Value = "3.214.96"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").Value = "5.93"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").Value = "0.519"
$ws.Range("E12").Value = "  +3.95%  "
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").Value = "39.48"
$ws.Range("E14").Value = "  +6.47%  "
$ws.Range("D15").Value = "3.743.06"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "66.963.56"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "7.52"
$ws.Range("E17").Value = "  +5.41%  "
$ws.Range("D18").Value = "3.218.64"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.112"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "525.30"
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("D21").Value = "15.55"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").Value = "0.748"
$ws.Range("E22").Value = "  +4.77%  "
$ws.Range("D23").Value = "8.24"
$ws.Range("E23").Value = "  +6.83%  "
$ws.Range("D24").Value = "15.15"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "85.64"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "9.36"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("D29").Value = "2.42"
$ws.Range("E29").Value = "  +10.92%  "
$ws.Range("D30").Value = "3.02"
$ws.Range("E30").Value = "  +8.59%  "
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  +11.00%  "
$ws.Range("D32").Value = "28.46"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").Value = "526.36"
$ws.Range("E36").Value = "  +11.06%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").Value = "0.0913"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("D39").Value = "0.0429"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("E40").Value = "  +10.00%  "
$ws.Range("D41").Value = "8.94"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").Value = "2.93"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "0.0₃0694"
$ws.Range("E43").Value = "  +17.49%  "
$ws.Range("D44").Value = "0.303"
$ws.Range("E44").Value = "  +7.60%  "
$ws.Range("D46").Value = "2.917.14"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").Value = "28.82"
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  +10.58%  "
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("D50").Value = "2.36"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("E51").Value = "  +0.00%  "
